$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.200.44"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").Value = "3.559.19"
$ws.Range("E3").Value = "  +2.01%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "621.00"
$ws.Range("E5").Value = "  +2.39%  "
$ws.Range("D6").Value = "155.08"
$ws.Range("E6").Value = "  +4.46%  "
$ws.Range("D7").Value = "3.556.27"
$ws.Range("E7").Value = "  +1.94%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +2.01%  "
$ws.Range("E10").Value = "  +5.58%  "
$ws.Range("D11").Value = "7.35"
$ws.Range("E11").Value = "  +5.62%  "
$ws.Range("D12").Value = "0.438"
$ws.Range("E12").Value = "  +3.89%  "
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").Value = "0.0000222"
$ws.Range("E13").Value = "  +1.96%  "
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "33.15"
$ws.Range("E14").Value = "  +5.49%  "
$ws.Range("D15").Value = "4.160.01"
$ws.Range("E15").Value = "  +1.98%  "
$ws.Range("D16").Value = "3.558.74"
$ws.Range("E16").Value = "  +2.16%  "
$ws.Range("D17").Value = "68.120.16"
$ws.Range("E17").Value = "  +1.36%  "
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("E19").Value = "  +6.44%  "
$ws.Range("D20").Value = "15.98"
$ws.Range("E20").Value = "  +6.27%  "
$ws.Range("D21").Value = "9.97"
$ws.Range("E21").Value = "  +10.53%  "
$ws.Range("D22").Value = "455.21"
$ws.Range("E22").Value = "  +1.83%  "
$ws.Range("D23").Value = "0.643"
$ws.Range("E23").Value = "  +3.65%  "
$ws.Range("D24").Value = "78.47"
$ws.Range("E24").Value = "  +1.71%  "
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").Value = "0.0000129"
$ws.Range("E25").Value = "  +1.48%  "
$ws.Range("D26").Value = "3.698.41"
$ws.Range("E26").Value = "  +1.94%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Value = "10.54"
$ws.Range("E28").Value = "  +4.16%  "
$ws.Range("D29").Value = "9.12"
$ws.Range("E29").Value = "  +10.06%  "
$ws.Range("D30").Value = "2.57"
$ws.Range("E30").Value = "  +3.12%  "
$ws.Range("D31").Value = "1.69"
$ws.Range("E31").Value = "  +8.89%  "
$ws.Range("D32").Value = "0.172"
$ws.Range("E32").Value = "  +6.05%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("E34").Value = "  +4.36%  "
$ws.Range("D35").Value = "26.10"
$ws.Range("E35").Value = "  +1.71%  "
$ws.Range("E36").Value = "  +3.65%  "
$ws.Range("D37").Value = "3.551.60"
$ws.Range("E37").Value = "  +2.06%  "
$ws.Range("D38").Value = "8.27"
$ws.Range("E38").Value = "  +3.55%  "
$ws.Range("E39").Value = "  +7.68%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("D41").Value = "178.72"
$ws.Range("E41").Value = "  +2.88%  "
$ws.Range("E42").Value = "  +5.44%  "
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("E44").Value = "  +3.68%  "
$ws.Range("D45").Value = "31.13"
$ws.Range("E45").Value = "  +16.09%  "
$ws.Range("D46").Value = "0.897"
$ws.Range("E46").Value = "  +1.72%  "
$ws.Range("D47").Value = "46.64"
$ws.Range("E47").Value = "  +2.67%  "
$ws.Range("D48").Value = "1.33"
$ws.Range("E48").Value = "  +6.82%  "
$ws.Range("D49").Value = "2.66"
$ws.Range("E49").Value = "  +4.02%  "
$ws.Range("D50").Value = "7.79"
$ws.Range("E50").Value = "  +3.52%  "
$ws.Range("D51").Value = "0.262"
$ws.Range("E51").Value = "  +7.39%  "
